# "Generate Report for Archive"
#
# 1) The status text "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4 - all backed by the
#    same shared string).
# 2) The "status" column(s) get narrower: width 17.2159881591797 ->
#    13.4101848602295 (Overview columns E & F; zh-cn and de-de column C).

$wb = $excel.ActiveWorkbook

# --- 1) Swap the status text on every sheet that uses it ---------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation", -4142, 1, $false, $false, $true) | Out-Null
}

# --- 2) Narrow the status column(s) -------------------------------------
# ColumnWidth of 12.5 is the Excel character-width value that lands on the
# same rendered column width as the target OOXML `width` (the nearest
# pixel-snapped width to 13.4101848602295).
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Columns.Item(5).ColumnWidth = 12.5   # column E ("zh-cn" status)
$ovw.Columns.Item(6).ColumnWidth = 12.5   # column F ("de-de" status)

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5  # column C ("Status")

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5  # column C ("Status")
